$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "67.928.43"
$ws.Range("E2").Value = "  +1.75%  "
Set-TextValue $ws "D3" "3.336.79"
$ws.Range("E3").Value = "  +1.74%  "
Set-TextValue $ws "D4" "1.00"
$ws.Range("E4").Value = "  +0.14%  "
Set-TextValue $ws "D5" "580.52"
$ws.Range("E5").Value = "  +1.66%  "
Set-TextValue $ws "D6" "177.25"
$ws.Range("E6").Value = "  +0.60%  "
Set-TextValue $ws "D7" "1.00"
$ws.Range("E7").Value = "  +0.08%  "
Set-TextValue $ws "D8" "0.590"
$ws.Range("E8").Value = "  +1.90%  "
Set-TextValue $ws "D9" "3.333.87"
$ws.Range("E9").Value = "  +1.85%  "
$ws.Range("E10").Value = "  +4.90%  "
$ws.Range("E11").Value = "  +1.62%  "
Set-TextValue $ws "D12" "46.85"
$ws.Range("E12").Value = "  +2.66%  "
$ws.Range("E13").Value = "  +1.66%  "
Set-TextValue $ws "D14" "686.84"
$ws.Range("E14").Value = "  -1.39%  "
Set-TextValue $ws "D15" "3.877.00"
$ws.Range("E15").Value = "  +1.91%  "
Set-TextValue $ws "D16" "8.45"
$ws.Range("E16").Value = "  +1.65%  "
Set-TextValue $ws "D17" "67.985.71"
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D18" "3.344.11"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws "D19" "0.118"
$ws.Range("E19").Value = "  -0.50%  "
Set-TextValue $ws "D20" "17.42"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("E21").Value = "  +3.30%  "
Set-TextValue $ws "D22" "0.898"
$ws.Range("E22").Value = "  +1.09%  "
Set-TextValue $ws "D23" "17.03"
$ws.Range("E23").Value = "  +0.72%  "
Set-TextValue $ws "D24" "5.34"
$ws.Range("E24").Value = "  +3.94%  "
Set-TextValue $ws "D25" "98.89"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E28").Value = "  +2.03%  "
Set-TextValue $ws "D29" "32.89"
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("E30").Value = "  +1.72%  "
Set-TextValue $ws "D31" "7.09"
$ws.Range("E31").Value = "  +4.69%  "
Set-TextValue $ws "D32" "575.17"
$ws.Range("E32").Value = "  +1.64%  "
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("E34").Value = "  +1.93%  "
Set-TextValue $ws "D36" "3.710.33"
$ws.Range("E36").Value = "  -4.78%  "
Set-TextValue $ws "D37" "56.98"
$ws.Range("E37").Value = "  +2.60%  "
Set-TextValue $ws "D38" "3.28"
$ws.Range("E38").Value = "  -0.95%  "
Set-TextValue $ws "D39" "34.67"
$ws.Range("E39").Value = "  +8.75%  "
$ws.Range("E40").Value = "  +2.64%  "
Set-TextValue $ws "D41" "2.67"
$ws.Range("E41").Value = "  +2.21%  "
Set-TextValue $ws "D42" "3.19"
$ws.Range("E42").Value = "  +6.20%  "
Set-TextValue $ws "D43" "3.38"
$ws.Range("E43").Value = "  +0.79%  "
Set-TextValue $ws "D44" "0.0₃0676"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("E46").Value = "  +0.23%  "
Set-TextValue $ws "D47" "2.67"
$ws.Range("E47").Value = "  +4.93%  "
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("E50").Value = "  -2.61%  "
Set-TextValue $ws "D51" "129.48"
$ws.Range("E51").Value = "  -0.38%  "
